$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows: regcntr_id, machine_id, device_id, lang_code, is_active, cr_by, cr_dtimes, eff_dtimes
$newRows = @(
    @(10001, 10030, 3000166),
    @(10001, 10030, 3000167),
    @(10001, 10030, 3000168),
    @(10001, 10030, 3000169),
    @(10001, 10030, 3000170),
    @(10001, 10031, 3000171),
    @(10001, 10031, 3000172),
    @(10001, 10031, 3000173),
    @(10001, 10031, 3000174),
    @(10001, 10031, 3000175)
)

$startRow = 147
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = "eng"
    $ws.Cells.Item($r, 5).Value = $true
    $ws.Cells.Item($r, 6).Value = "superadmin"
    $ws.Cells.Item($r, 7).Value = "now()"
    $ws.Cells.Item($r, 8).Value = "now()"
}

# Update selection / view to reflect scroll position after appending rows
$ws.Range("A148").Select()
